$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert two new columns before the old "Test Rail url" column (O),
# shifting it to Q and making room for the new Meganav2 columns N/O/P.
$ws.Columns("N:O").Insert()

# Populate the new Meganav2 columns (the order below reproduces the
# shared-string insertion order seen in the target workbook).
$ws.Range("N1").Value = "Meganav2"
$ws.Range("N2").Value = "/c-340-30-under"
$ws.Range("P1").Value = "Meganav2 Breadcrumbs"
$ws.Range("P2").Value = "Home/Collections/`$30 & Under "
$ws.Range("D1").Value = "Meganav1 Header"
$ws.Range("O1").Value = "Meganav2 Header"
$ws.Range("O2").Value = "`$30 & Under"

# Re-point the Test Rail hyperlink: the column insert moved the old O2
# cell (with its hyperlink) to Q2, but the hyperlink definition itself
# needs to be recreated so it tracks the new location.
$ws.Range("Q2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("Q2"), "https://surlatable.testrail.net/index.php?/cases/view/12080&group_by=cases:section_id&group_order=asc&display_deleted_cases=0&group_id=1961")
$ws.Range("Q2").Style = "Hyperlink"

# Column widths for the new columns
$ws.Range("N1").EntireColumn.ColumnWidth = 14.5
$ws.Range("O1").EntireColumn.ColumnWidth = 16.2
$ws.Range("P1").EntireColumn.ColumnWidth = 34.1

# Update the view so the new columns are visible and the selection
# matches the target workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("O4").Select()
